# Added Week 15 simulations
$wb = $excel.ActiveWorkbook

# --- OFF sheet (row 2) ---
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B2").Value = 167
$wsOff.Range("C2").Value = 114
$wsOff.Range("D2").Value = 46
$wsOff.Range("E2").Value = 22
$wsOff.Range("F2").Value = 3
$wsOff.Range("G2").Value = 3

# --- DEF sheet (row 2) ---
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B2").Value = 171
$wsDef.Range("C2").Value = 117
$wsDef.Range("D2").Value = 45
$wsDef.Range("E2").Value = 25
$wsDef.Range("F2").Value = 3
